$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange
$lastRow = $ur.Rows.Count + $ur.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text
    if ($val -ne $null -and $val.Length -gt 0 -and $val.Contains(",")) {
        $parts = $val -split ","
        $trimmed = @()
        foreach ($p in $parts) { $trimmed += $p.Trim() }

        if ($trimmed -contains "System") {
            $reversed = @()
            for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
                $reversed += $trimmed[$i]
            }
            $newVal = [string]::Join(", ", $reversed)
            $cell.Value2 = $newVal
        }
    }
}
